$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column A (everything else shifts right: A:E -> B:F)
$ws.Columns.Item(1).Insert()

# New "ID" column: header + first (only) record's id
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = 1

# Widen the (now) date column (D) to fit its contents better
$ws.Columns.Item(4).ColumnWidth = 11.92

# Switch the workbook's base font from Calibri to Arial
$wb.Styles("Normal").Font.Name = "Arial"

# Changing the Normal style can disturb explicit number formats on existing
# cells in this runtime - restore the date cell's display format.
$ws.Range("D2").NumberFormat = "d-mmm"

# Final selected cell observed after the edits
$ws.Range("A4").Select()
